$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content edits --------------------------------------------------
# A2: "footsteps on wood" -> "footsteps on pavement"
$ws.Range("A2").Value2 = "footsteps on pavement"

# New column D (header "voices" + 6 values below it)
$ws.Range("D1").Value2 = "voices"
$ws.Range("D2").Value2 = "brother"
$ws.Range("D3").Value2 = "sailor pete"
$ws.Range("D4").Value2 = "bully1"
$ws.Range("D5").Value2 = "bully2"
$ws.Range("D6").Value2 = "vendor"
$ws.Range("D7").Value2 = "vendors daughter"

# New row 7 (A7 / B7)
$ws.Range("A7").Value2 = "seagulls"
$ws.Range("B7").Value2 = "rock being smacked"

# Give the new header cell D1 the same formatting (bold font + shaded
# fill) as the other header cells by copying the format from C1.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Column widths ---------------------------------------------------------
# Column A widens to fit the longer "footsteps on pavement" text;
# column D is the new column and needs its own width.
$ws.Columns.Item(1).ColumnWidth = 21.166666666666668
$ws.Columns.Item(4).ColumnWidth = 15.877604166666666

# --- Selection ---------------------------------------------------------
$ws.Range("B8").Select() | Out-Null
